$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.161.62"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").Value = "3.505.50"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.19%  "
$ws.Range("E7").Value = "  -1.55%  "
$ws.Range("D8").Value = "3.499.18"
$ws.Range("E8").Value = "  -1.21%  "
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.25"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.585"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.12"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.31%  "
$ws.Range("E14").Value = "  -1.89%  "
$ws.Range("D15").Value = "4.071.50"
$ws.Range("E15").Value = "  -1.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.35"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "612.12"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.81%  "
$ws.Range("D18").Value = "3.502.80"
$ws.Range("E18").Value = "  -1.33%  "
$ws.Range("D19").Value = "70.113.08"
$ws.Range("E19").Value = "  -1.07%  "
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("E22").Value = "  -1.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "98.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "15.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.41%  "
$ws.Range("E26").Value = "  -3.80%  "
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("E28").Value = "  -2.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.27%  "
$ws.Range("E31").Value = "  -4.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.25%  "
$ws.Range("E33").Value = "  -5.02%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.51%  "
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "628.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0995"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.85%  "
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("E38").Value = "  +6.19%  "
$ws.Range("E39").Value = "  -4.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "56.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.24%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("D43").Value = "3.364.34"
$ws.Range("E43").Value = "  +0.56%  "
$ws.Range("D44").Value = "0.0₃0732"
$ws.Range("E44").Value = "  +1.62%  "
$ws.Range("E45").Value = "  -6.01%  "
$ws.Range("E46").Value = "  -4.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "31.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.93%  "
$ws.Range("E48").Value = "  -4.24%  "
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.15%  "
$ws.Range("E51").Value = "  -0.03%  "
